# Applies the "Debate 1 -> Debate 3" revision:
#  1. Title text: "Debate 1" -> "Debate 3"
#  2. Date text:  "F2025" -> "S2026"
#  3. New "Overview" / Heading2 paragraph inserted after the Date paragraph
#  4. FirstParagraph text: "Debate stuff" -> "Nothing to see here yet"
#  5. A document bookmark named "overview" wrapping the new Heading2
#     paragraph and the FirstParagraph paragraph that follows it.

$d = $word.ActiveDocument

# 1. Title
$d.Content.Find.Execute("Debate 1", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Debate 3", 2) | Out-Null

# 2. Date
$d.Content.Find.Execute("F2025", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "S2026", 2) | Out-Null

# 3. Locate the Date paragraph and insert a brand-new, empty paragraph right
#    after it; that new paragraph becomes the "Overview" heading.
$datePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Date") {
        $datePara = $p
        break
    }
}
$datePara.Range.InsertParagraphAfter() | Out-Null

$overviewPara = $null
$foundDate = $false
foreach ($p in $d.Paragraphs) {
    if ($foundDate) {
        $overviewPara = $p
        break
    }
    if ($p.Style.NameLocal -eq "Date") {
        $foundDate = $true
    }
}
$overviewPara.Style = "Heading2"
$overviewPara.Range.Text = "Overview"

# 4. FirstParagraph text
$d.Content.Find.Execute("Debate stuff", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Nothing to see here yet", 2) | Out-Null

# 5. Wrap the new Heading2 paragraph and the FirstParagraph paragraph that
#    follows it in a bookmark named "overview".
$overviewPara = $null
$followingPara = $null
foreach ($p in $d.Paragraphs) {
    if ($overviewPara -ne $null -and $followingPara -eq $null) {
        $followingPara = $p
    }
    if ($p.Style.NameLocal -eq "Heading 2" -and $p.Range.Text.TrimEnd() -eq "Overview") {
        $overviewPara = $p
    }
}
$bkRange = $d.Range($overviewPara.Range.Start, $followingPara.Range.End)
$d.Bookmarks.Add("overview", $bkRange) | Out-Null
